# Add two new columns, I ("I0") and J ("IF"), to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the rest of row 1 (bold font,
# thin border, centered/top aligned) by copying H1's format over.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (rows 2-41) ---
$data = @(
    @(2, 9, 9),
    @(3, 5, 6),
    @(4, 9, 9),
    @(5, 5, 6),
    @(6, 8, 8),
    @(7, 7, 7),
    @(8, 6, 6),
    @(9, 9, 9),
    @(10, 9, 9),
    @(11, 7, 8),
    @(12, 7, 7),
    @(13, 8, 8),
    @(14, 5, 5),
    @(15, 6, 7),
    @(16, 7, 7),
    @(17, 8, 8),
    @(18, 7, 7),
    @(19, 7, 7),
    @(20, 7, 8),
    @(21, 9, 9),
    @(22, 7, 8),
    @(23, 8, 8),
    @(24, 9, 9),
    @(25, 6, 7),
    @(26, 3, 5),
    @(27, 8, 8),
    @(28, 5, 5),
    @(29, 8, 8),
    @(30, 7, 7),
    @(31, 8, 8),
    @(32, 7, 8),
    @(33, 5, 6),
    @(34, 6, 7),
    @(35, 8, 8),
    @(36, 8, 9),
    @(37, 9, 9),
    @(38, 8, 8),
    @(39, 7, 7),
    @(40, 3, 3),
    @(41, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $i0 = $row[1]
    $if = $row[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $if
}

Write-Output "I0 and IF columns added"
